# Mark checklist items "1." and "2." as done by applying strikethrough
# formatting to their entire paragraphs (runs + paragraph mark), matching
# the treatment already present on items 4, 5, and 7 in this document.

$d = $word.ActiveDocument

$item1 = $d.Paragraphs.Item(3)
$item1.Range.Font.StrikeThrough = 1

$item2 = $d.Paragraphs.Item(4)
$item2.Range.Font.StrikeThrough = 1
